$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("F1_Venta_23_Ene_Porcentaje") was entered as whole numbers
# (101, 51, ...) but really represents a percentage; convert the stored
# values to their fractional form (1.01, 0.51, ...) and display them with
# a one-decimal percentage number format.
$rng = $ws.Range("C2:C13")

$values = @(1.01, 0.51, 0.67, 0.47, 0.87, 0.82, 0.73, 0.68, 0.86, 0.71, 0.66, 0.42)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}

$rng.NumberFormat = "0.0%"
$rng.WrapText = $true
$rng.VerticalAlignment = -4108
$rng.Font.ThemeColor = 1

# Reflect where the user last clicked after making the edit.
$null = $ws.Range("C17").Select()
